$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 25-26, shifting existing rows down
$ws.Rows("25:26").Insert()

# Fill in the data for the two newly inserted rows (25 and 26)
# Row 25
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44565
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103002
$ws.Range("J25").Value = "Ciruela"
$ws.Range("K25").Value = "Black Amber"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("Q25").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 778
$ws.Range("T25").Value = 18

# Row 26
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44565
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100103
$ws.Range("H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I26").Value = 100103002
$ws.Range("J26").Value = "Ciruela"
$ws.Range("K26").Value = "Black Amber"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("Q26").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 667
$ws.Range("T26").Value = 18
